# Applies the "adds parsing of Regions" edit:
#  - settings sheet: append two new rows (outfile / contextfile)
#  - regions sheet: rename header I1 "enableatstart" -> "enabledatstart",
#    and add theme/level (M/N) values for the NOGPS and BG rows
#  - switch the active sheet/tab from "settings" to "regions" and update
#    the remembered selection on each sheet

$wb = $excel.ActiveWorkbook

$settings = $wb.Worksheets.Item("settings")
$regions  = $wb.Worksheets.Item("regions")

# --- regions sheet: fix column header typo ---
$regions.Range("I1").Value = "enabledatstart"

# --- settings sheet: two new rows ---
$settings.Range("A7").Value = "outfile"
$settings.Range("B7").Value = "daoplayer-test1.json"
$settings.Range("A8").Value = "contextfile"
$settings.Range("B8").Value = "context-test1.json"

# --- regions sheet: add parsed theme/level for the NOGPS and BG rows ---
$regions.Range("M2").Value = "t"
$regions.Range("N2").Value = 1

$regions.Range("M3").Value = "t2"
$regions.Range("N3").Value = 1

# --- update selections left on each sheet ---
[void]$settings.Range("A9").Select()
[void]$regions.Range("N3").Select()

# --- switch active sheet to "regions" ---
[void]$regions.Activate()
